$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-4, columns E:T
$data = @{
    2 = @(2, 0.6666666666666666, 0.1823923333333334, 0.547177, 1, 1, 3, 1, 3.083576666666666, 9.250729999999999, 0.2272509363535097, 0.2272509363535097, 0.5624207432455556, 5.06178668921, 0.2272509363535097, 0.2272509363535097)
    3 = @(2, 0.6666666666666666, 0.1823923333333334, 0.547177, 1, 1, 3, 1, 6.453984666666667, 19.361954, 0.4756405360586227, 0.4756405360586227, 1.177157322650889, 10.594415903858, 0.4756405360586227, 0.4756405360586227)
    4 = @(2, 0.6666666666666666, 0.1823923333333334, 0.547177, 1, 1, 3, 1, 4.031477000000001, 12.094431, 0.2971085275878677, 0.2971085275878677, 0.7353104968096669, 6.617794471287001, 0.2971085275878677, 0.2971085275878677)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i   # Column E is index 5
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
